$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Data Test for web/mobile/api": the positive-flow API test case (row 10)
# now uses a fresh throw-away email address. Update D10 in place so the
# existing hyperlink relationship (mailto:) is kept untouched and only the
# displayed/stored text changes from archen17@gmail.com to archen22@gmail.com.
$ws.Range("D10").Value = "archen22@gmail.com"

# The author had scrolled the sheet down a bit and had cell D10 selected
# when they saved (previously F8 was selected at the very top of the
# sheet). Reflect the same view/selection state.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D10").Select()
